$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Append 10 new reconciliation rows (57-66) ---
#
# Columns B (Noms), C (Routes) and D (Sous-Zone)/I (Site) mix brand new
# text values with values that already exist elsewhere in the sheet.
# We deliberately write column-by-column (all of B, then all of C, then
# all of D, then all of I) instead of row-by-row so that the *new*
# strings are introduced to the workbook in the same grouped order as
# the source data (all 10 new names together, then the new route code,
# then the new site name) rather than interleaved.

# Numero (A)
$ws.Cells.Item(57, 1).Value = 237654101067
$ws.Cells.Item(58, 1).Value = 237675453374
$ws.Cells.Item(59, 1).Value = 237680039383
$ws.Cells.Item(60, 1).Value = 237678973363
$ws.Cells.Item(61, 1).Value = 237681663743
$ws.Cells.Item(62, 1).Value = 237651646213
$ws.Cells.Item(63, 1).Value = 237652194260
$ws.Cells.Item(64, 1).Value = 237671615641
$ws.Cells.Item(65, 1).Value = 237653816480
$ws.Cells.Item(66, 1).Value = 237673593310

# Noms (B)
$ws.Cells.Item(57, 2).Value = "MAKUETCHE TCHEHGHIE CELINE GIRESSE CHIC MOBILE SARL"
$ws.Cells.Item(58, 2).Value = "ABEL MOUNTAPMBEME"
$ws.Cells.Item(59, 2).Value = "SPECTRUM LTDLA CBOX R0 CEDRICK MARCIALLE WANDJI"
$ws.Cells.Item(60, 2).Value = "MOSSU TAGNE ANNE FLORE TOP MOBIL"
$ws.Cells.Item(61, 2).Value = "LA NEGRESSE SARL FONGA SINTCHA YOLANDE MIREILLE"
$ws.Cells.Item(62, 2).Value = "AMADOU AHIJO ETS MOBILE FINANCIAL SERVICES MFS"
$ws.Cells.Item(63, 2).Value = "CRISTELLE DIANE TCHAHANE"
$ws.Cells.Item(64, 2).Value = "BEGO FOGUE CHRISTELLE KAMILAH CONNECTION GROUP"
$ws.Cells.Item(65, 2).Value = "BERYL NAKOMA TOUFOIN TOP MOBIL TELECOM"
$ws.Cells.Item(66, 2).Value = "JEAN JACQUES YENDJE"

# Routes (C)
$ws.Cells.Item(57, 3).Value = "Rte_4"
$ws.Cells.Item(58, 3).Value = 0
$ws.Cells.Item(59, 3).Value = 0
$ws.Cells.Item(60, 3).Value = "Rte_8"
$ws.Cells.Item(61, 3).Value = 0
$ws.Cells.Item(62, 3).Value = 0
$ws.Cells.Item(63, 3).Value = "Rte_6"
$ws.Cells.Item(64, 3).Value = 0
$ws.Cells.Item(65, 3).Value = "Rte_3"
$ws.Cells.Item(66, 3).Value = 0

# Sous-Zone (D)
$ws.Cells.Item(57, 4).Value = "Esg Building"
$ws.Cells.Item(58, 4).Value = "Makepe Conquete"
$ws.Cells.Item(59, 4).Value = "Agape Ocm"
$ws.Cells.Item(60, 4).Value = "Total Ndokotti"
$ws.Cells.Item(61, 4).Value = "Cite Bassa"
$ws.Cells.Item(62, 4).Value = "Ndogbong Vallee Ocm"
$ws.Cells.Item(63, 4).Value = "Ndogbong-Citadelle"
$ws.Cells.Item(64, 4).Value = "Makepe Conquete"
$ws.Cells.Item(65, 4).Value = "Mobil Guiness"
$ws.Cells.Item(66, 4).Value = "Ndokoti Carrefour"

# Montants OOS (E)
$ws.Cells.Item(57, 5).Value = 35000
$ws.Cells.Item(58, 5).Value = 44030
$ws.Cells.Item(59, 5).Value = 10000
$ws.Cells.Item(60, 5).Value = 163510
$ws.Cells.Item(61, 5).Value = 128170
$ws.Cells.Item(62, 5).Value = 16509.8
$ws.Cells.Item(63, 5).Value = 112790
$ws.Cells.Item(64, 5).Value = 10510
$ws.Cells.Item(65, 5).Value = 344600
$ws.Cells.Item(66, 5).Value = 61380

# Balance (F)
$ws.Cells.Item(57, 6).Value = 82435
$ws.Cells.Item(58, 6).Value = 333287
$ws.Cells.Item(59, 6).Value = 15241
$ws.Cells.Item(60, 6).Value = 217587
$ws.Cells.Item(61, 6).Value = 335141
$ws.Cells.Item(62, 6).Value = 4643
$ws.Cells.Item(63, 6).Value = 436173
$ws.Cells.Item(64, 6).Value = 7167
$ws.Cells.Item(65, 6).Value = 1929172
$ws.Cells.Item(66, 6).Value = 123365

# Valeur Calculee (G)
$ws.Cells.Item(57, 7).Value = 47435
$ws.Cells.Item(58, 7).Value = 289257
$ws.Cells.Item(59, 7).Value = 5241
$ws.Cells.Item(60, 7).Value = 54077
$ws.Cells.Item(61, 7).Value = 206971
$ws.Cells.Item(62, 7).Value = -11866.8
$ws.Cells.Item(63, 7).Value = 323383
$ws.Cells.Item(64, 7).Value = -3343
$ws.Cells.Item(65, 7).Value = 1584572
$ws.Cells.Item(66, 7).Value = 61985

# Jours de Stock (H)
$ws.Cells.Item(57, 8).Value = 2.355285714285714
$ws.Cells.Item(58, 8).Value = 7.569543493072905
$ws.Cells.Item(59, 8).Value = 1.5241
$ws.Cells.Item(60, 8).Value = 1.330725949483212
$ws.Cells.Item(61, 8).Value = 2.614816259655146
$ws.Cells.Item(62, 8).Value = 0.2812269076548474
$ws.Cells.Item(63, 8).Value = 3.867124745101516
$ws.Cells.Item(64, 8).Value = 0.6819219790675547
$ws.Cells.Item(65, 8).Value = 5.598293673824724
$ws.Cells.Item(66, 8).Value = 2.009856630824373

# Site (I)
$ws.Cells.Item(57, 9).Value = "Ndogbong"
$ws.Cells.Item(58, 9).Value = "Ndogbong"
$ws.Cells.Item(59, 9).Value = "Cite Sic"
$ws.Cells.Item(60, 9).Value = "Ndogbong"
$ws.Cells.Item(61, 9).Value = "Cite Sic"
$ws.Cells.Item(62, 9).Value = "Ndogbong"
$ws.Cells.Item(63, 9).Value = "Ndogbong"
$ws.Cells.Item(64, 9).Value = "Ndogbong"
$ws.Cells.Item(65, 9).Value = "Cite Sic"
$ws.Cells.Item(66, 9).Value = "Ndogbong"
